$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:205 down to 104:206
$ws.Rows("103:103").Insert()

# Populate the newly inserted row 103 with the new data point
$ws.Range("A103").Value = 1
$ws.Range("B103").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C103").Value = "Arica y Parinacota"
$ws.Range("D103").Value = 44629
$ws.Range("E103").Value = 15
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100108
$ws.Range("H103").Value = "Tropicales y subtropicales"
$ws.Range("I103").Value = 100108006
$ws.Range("J103").Value = "Plátano"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Pintón"
$ws.Range("M103").Value = 120
$ws.Range("N103").Value = 16000
$ws.Range("O103").Value = 17000
$ws.Range("P103").Value = 16500
$ws.Range("Q103").Value = "$/caja 20 kilos"
$ws.Range("R103").Value = "Bolivia"
$ws.Range("S103").Value = 825
$ws.Range("T103").Value = 20
